$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Minimize the workbook window (bookViews -> workbookView minimized="1")
$excel.WindowState = -4140  # xlMinimized
$wb.Windows.Item(1).WindowState = -4140

# Rename the field "Operator_maszyny" -> "Nazwa_op_maszyny" everywhere it is
# used as a column/field-name header: in the "Maszyny" table (A13), the
# "Operator_maszyny_Maszyny" table (C8) and the "Operator_Maszyny" table (E6).
# (The E4 cell, which uses the same original text as a section title, is left
# untouched - it is not part of this rename in the source edit.)
$ws.Range("A13").Value = "Nazwa_op_maszyny"
$ws.Range("C8").Value = "Nazwa_op_maszyny"
$ws.Range("E6").Value = "Nazwa_op_maszyny"

# Move/record the active selection on the sheet to C8
$ws.Range("C8").Select()
